$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: new article (Ser 23, dated 23-Jan-2020)
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = 43853
$ws.Range("B24").NumberFormat = "d-mmm-yy"

# Set the Ayats (C) and Tags (F) strings before Content (D) so the new
# shared-string table entries line up the same way the source workbook has them.
$ws.Range("C24").Value = 'Surah Baqarah, 236 - 249'
$ws.Range("F24").Value = 'Importance of prayers, Sadaqah, Charity, Financial Comfort, Tips to succeed'
$ws.Range("D24").Value = 'h3: Tips to live a problem free life
p.note: These articles are happening, because I managed to say 5 prayers a day consistently for past 40 days. Prayers help us live a productive life. 
quote: Maintain with care the [obligatory] prayers and [in particular] the middle prayer and stand before Allah, devoutly obedient. And if you fear [an enemy, then pray] on foot or riding. But when you are secure, then remember Allah [in prayer], as He has taught you that which you did not [previously] know. <br> - Surah Baqarah verse 238, 239
p: Prayers are important, emphasised a number of times by our parents, elders, colleagues, most muslim scholars and now in Quran. 
h3: Protecting our prayers
p: Below are some tips that have helped me in protecting my prayers:-
p.b-left: <b>1. Draw a picture </b>of what we want to achieve in our life and hang it somewhere prominent.
p.b-left: <b>2. Live in </b>your dream picture. Feel it happening around you.
p.b-left: <b>3. Print a 30 Days Challenge calendar</b> and hang it on the wall, somewhere prominent. In your wardrobe, wall or drawing room.
p.b-left: <b>4. Say your 5 prayers </b>a day aggressively with Jamat in your local mosque. Do not miss even a single prayer for next 30 days. Put a cross on it, every night before you sleep.
p.note: If you miss jamat due to some serious problem, it is alright. Say your prayer later and keep it accounted for.
p.b-left: <b>5. Throw yourself in air</b>, here and there during this period. Make it more <b>risky</b> and fun. 
p: I am writing articles everyday, I do not know the entire day what I am going to write. It is all risk. My heart beats faster and I am more on my toes. Everyday I need more favors from my Allah for these articles to work. Adrenaline rush keeps my heart in my mouth round the clock. 
h3: Are you undergoing a financial crises?
quote: Who is it that would loan Allah a goodly loan so He may multiply it for him many times over? And it is Allah who withholds and grants abundance, and to Him you will be returned. <br> - Surah Baqarah verse 245
p: Tips to get rid of financial burden:-
p.b-left: <b>1. Buy a box </b>, similar to the one below.
p.b-left: <b>2. Place it </b>at your’s home entrance.
p.b-left: <b>3. Every morning</b> put some sadqah in it before leaving the house.
p.b-left: <b>4. Make somebody </b>richer every month.
img.width-half: box.jpg
p: Do not miss out the habit of giving sadaqah every day. It is one fun way to protect ourselves from evils outside and maintain our financial well-being. 
p.note: Guys, if you like this project. Please follow this project''s page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>.'
$ws.Range("E24").Value = "Qasim Ali"

$ws.Rows.Item(24).RowHeight = 409.6

[void]$ws.Range("D24").Select()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
